# Natmi following Dr Hou advice
# Recomputed LR-pair table for Pdgfc-Pdgfra: "ECs" is now included as a
# Sending cluster (in addition to already being a Target cluster), growing
# the 3x4 cluster grid to a full 4x4 grid (rows 2-17), and refreshed
# expression/specificity statistics for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Pdgfc", "Pdgfra", "ECs", 2, 0.6666666666666666, 0.276989, 0.830967, 0.0126770351133645, 0.0126770351133645, 3, 1, 0.236283, 0.708849, 0.001461516295904947, 0.001461516295904947, 0.06544779188699999, 0.589030126983, 0.00001852769340194144, 0.00001852769340194144)
    ,@("ECs", "Pdgfc", "Pdgfra", "FAPs", 2, 0.6666666666666666, 0.276989, 0.830967, 0.0126770351133645, 0.0126770351133645, 3, 1, 160.8390196666666, 482.517059, 0.9948614511421033, 0.9948614511421032, 44.55063921845032, 400.955752966053, 0.01261189354906121, 0.01261189354906121)
    ,@("ECs", "Pdgfc", "Pdgfra", "M2", 2, 0.6666666666666666, 0.276989, 0.830967, 0.0126770351133645, 0.0126770351133645, 3, 1, 0.119972, 0.359916, 0.0007420806111836584, 0.0007420806111836582, 0.033230924308, 0.299078318772, 0.00000940738196492223, 0.000009407381964922224)
    ,@("ECs", "Pdgfc", "Pdgfra", "sCs", 2, 0.6666666666666666, 0.276989, 0.830967, 0.0126770351133645, 0.0126770351133645, 3, 1, 0.474493, 1.423479, 0.002934951950808252, 0.002934951950808251, 0.131429341577, 1.182864074193, 0.00003720648893643386, 0.00003720648893643384)
    ,@("FAPs", "Pdgfc", "Pdgfra", "ECs", 3, 1, 5.827140333333332, 17.481421, 0.2666924051719354, 0.2666924051719353, 3, 1, 0.236283, 0.708849, 0.001461516295904947, 0.001461516295904947, 1.376854199381, 12.391687794429, 0.0003897752961528684, 0.0003897752961528682)
    ,@("FAPs", "Pdgfc", "Pdgfra", "FAPs", 3, 1, 5.827140333333332, 17.481421, 0.2666924051719354, 0.2666924051719353, 3, 1, 160.8390196666666, 482.517059, 0.9948614511421033, 0.9948614511421032, 937.2315386734263, 8435.083848060836, 0.2653219932179294, 0.2653219932179293)
    ,@("FAPs", "Pdgfc", "Pdgfra", "M2", 3, 1, 5.827140333333332, 17.481421, 0.2666924051719354, 0.2666924051719353, 3, 1, 0.119972, 0.359916, 0.0007420806111836584, 0.0007420806111836582, 0.6990936800706666, 6.291843120635999, 0.0001979072630280297, 0.0001979072630280295)
    ,@("FAPs", "Pdgfc", "Pdgfra", "sCs", 3, 1, 5.827140333333332, 17.481421, 0.2666924051719354, 0.2666924051719353, 3, 1, 0.474493, 1.423479, 0.002934951950808252, 0.002934951950808251, 2.764937298184333, 24.884435683659, 0.0007827293948251164, 0.0007827293948251161)
    ,@("M2", "Pdgfc", "Pdgfra", "ECs", 3, 1, 2.750881, 8.252642999999999, 0.1259003607713204, 0.1259003607713204, 3, 1, 0.236283, 0.708849, 0.001461516295904947, 0.001461516295904947, 0.6499864153229999, 5.849877737906999, 0.0001840054289275967, 0.0001840054289275966)
    ,@("M2", "Pdgfc", "Pdgfra", "FAPs", 3, 1, 2.750881, 8.252642999999999, 0.1259003607713204, 0.1259003607713204, 3, 1, 160.8390196666666, 482.517059, 0.9948614511421033, 0.9948614511421032, 442.4490032596595, 3982.041029336936, 0.1252534156162702, 0.1252534156162701)
    ,@("M2", "Pdgfc", "Pdgfra", "M2", 3, 1, 2.750881, 8.252642999999999, 0.1259003607713204, 0.1259003607713204, 3, 1, 0.119972, 0.359916, 0.0007420806111836584, 0.0007420806111836582, 0.330028695332, 2.970258257988, 0.00009342821666942452, 0.00009342821666942447)
    ,@("M2", "Pdgfc", "Pdgfra", "sCs", 3, 1, 2.750881, 8.252642999999999, 0.1259003607713204, 0.1259003607713204, 3, 1, 0.474493, 1.423479, 0.002934951950808252, 0.002934951950808251, 1.305273778333, 11.747464004997, 0.0003695115094532495, 0.0003695115094532493)
    ,@("sCs", "Pdgfc", "Pdgfra", "ECs", 3, 1, 12.994657, 38.983971, 0.5947301989433799, 0.5947301989433798, 3, 1, 0.236283, 0.708849, 0.001461516295904947, 0.001461516295904947, 3.070416539931, 27.633748859379, 0.0008692078774225409, 0.0008692078774225405)
    ,@("sCs", "Pdgfc", "Pdgfra", "FAPs", 3, 1, 12.994657, 38.983971, 0.5947301989433799, 0.5947301989433798, 3, 1, 160.8390196666666, 482.517059, 0.9948614511421033, 0.9948614511421032, 2090.047892784587, 18810.43103506129, 0.5916741487588427, 0.5916741487588425)
    ,@("sCs", "Pdgfc", "Pdgfra", "M2", 3, 1, 12.994657, 38.983971, 0.5947301989433799, 0.5947301989433798, 3, 1, 0.119972, 0.359916, 0.0007420806111836584, 0.0007420806111836582, 1.558994989604, 14.030954906436, 0.000441337749521282, 0.0004413377495212818)
    ,@("sCs", "Pdgfc", "Pdgfra", "sCs", 3, 1, 12.994657, 38.983971, 0.5947301989433799, 0.5947301989433798, 3, 1, 0.474493, 1.423479, 0.002934951950808252, 0.002934951950808251, 6.165873783901, 55.492864055109, 0.001745504557593452, 0.001745504557593452)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}
